$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: account holder name
$ws.Range("C2").Value = "Hartmut"

# Account number (16-digit) must stay text, not be coerced to a number -
# format the cell as Text before assigning, then restore General formatting
# (lowercase "general" keeps the existing numeric-format id instead of
# minting a new custom format entry).
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("B3").NumberFormat = "general"

$ws.Range("C3").Value = "Mohaupt"

# Opening balance date
$ws.Range("D5").Value = "KONTOSTAND AM 11.07.2025"

# Row 6 - existing transaction updated
$ws.Range("B6").Value = "13.07."
$ws.Range("C6").Value = "14.07."
$ws.Range("D6").Value = "MCDONALDS Rothenburg ob der Tauber"
$ws.Range("E6").Value = "26,85-"

# Row 7 - existing transaction updated
$ws.Range("B7").Value = "17.07."
$ws.Range("C7").Value = "18.07."
$ws.Range("D7").Value = "AMAZON.DE MKTPLC EU AOXGDW"
$ws.Range("E7").Value = "162,72-"

# Row 8 - existing transaction updated
$ws.Range("B8").Value = "19.07."
$ws.Range("C8").Value = "20.07."
$ws.Range("D8").Value = "PAYPAL GOHZGP"
$ws.Range("E8").Value = "24,60-"

# Row 9 - existing transaction updated
$ws.Range("B9").Value = "23.07."
$ws.Range("C9").Value = "24.07."
$ws.Range("D9").Value = "EBAY MKTPLC EU TXXDZM"
$ws.Range("E9").Value = "152,27-"

# Row 10 - existing transaction updated
$ws.Range("B10").Value = "26.07."
$ws.Range("C10").Value = "27.07."
$ws.Range("D10").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 58494300"
$ws.Range("E10").Value = "84,50-"

# Row 11 - was an empty filler row, now a new transaction row. E11 needs to
# switch from the empty-row style to the right-aligned amount style used by
# the other transaction rows (drop the wrap-text / vertical-center that the
# blank-row style carried).
$ws.Range("B11").Value = "29.07."
$ws.Range("C11").Value = "30.07."
$ws.Range("D11").Value = "BURGER KING Parsberg"
$ws.Range("E11").Value = "28,00-"
$ws.Range("E11").VerticalAlignment = -4107
$ws.Range("E11").WrapText = $false

# Closing balance date/value
$ws.Range("D12").Value = "KONTOSTAND AM 02.08.2025"
$ws.Range("E12").Value = "478,94-"

# Next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 09.08.2025"
